# V 2.0.2 se arreglo la fechar y hora de reimpresion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Patient name / record number
$ws.Range("A6").Value = "GOMEZ  CHIVALAN  ELISA  NICOL"
$ws.Range("G6").Value = "/201762651"

# Birth date / age
$ws.Range("A9").Value = "2010-04-23"
$ws.Range("D9").Value = "6A"

# Nacionalidad (fix typo to match Lugar de Nacimiento)
$ws.Range("E11").Value = "GUATEMALA"

# Emergency contact
$ws.Range("A13").Value = "JUANA GOMEZ"
$ws.Range("D13").Value = "MAMA"
$ws.Range("E13").Value = "7MA C. 14-54 LA BARREDA Z.18"

# Fecha/hora de la asistencia médica + area de urgencia
$ws.Range("D14").Value = "Hora: 15:47:10"
$ws.Range("E14").Value = "Area de urgencia: MEDICINA"
$ws.Range("A15").Value = "24/10/2017"

# Tipo de Consulta
$ws.Range("D16").Value = "urgencia"
